$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from row 52 down to rows 53-64 to match date/text cell styles
$ws.Range("A52:AE52").Copy()
$ws.Range("A53:AE64").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 53
$ws.Cells.Item(53, 1).Value = 42531.559451747686
$ws.Cells.Item(53, 2).Value = "Bernardo Henz"
$ws.Cells.Item(53, 3).Value = 25.0
$ws.Cells.Item(53, 4).Value = "Masculino"
$ws.Cells.Item(53, 5).Value = "Pós-graduação (Mestrado/Doutorado/Pós-doc) COMPLETO"
$ws.Cells.Item(53, 6).Value = "Ciência da Computação"
$ws.Cells.Item(53, 7).Value = "Sempre direita"
$ws.Cells.Item(53, 8).Value = "Normalmente direita"
$ws.Cells.Item(53, 9).Value = "Sempre direita"
$ws.Cells.Item(53, 10).Value = "Sem prefêrencia"
$ws.Cells.Item(53, 11).Value = "Normalmente direita"
$ws.Cells.Item(53, 12).Value = "Normalmente esquerda"
$ws.Cells.Item(53, 13).Value = "Sempre direita"
$ws.Cells.Item(53, 14).Value = "mais que 12 horas"
$ws.Cells.Item(53, 15).Value = "mais que 12 horas"
$ws.Cells.Item(53, 16).Value = "Nunca"
$ws.Cells.Item(53, 17).Value = "Nunca"
$ws.Cells.Item(53, 18).Value = 5.0
$ws.Cells.Item(53, 19).Value = "mais que 12 horas"
$ws.Cells.Item(53, 20).Value = "até 4 horas"
$ws.Cells.Item(53, 21).Value = "Kinect, Wii, HMD"
$ws.Cells.Item(53, 22).Value = 3.0
$ws.Cells.Item(53, 23).Value = "Não"
$ws.Cells.Item(53, 24).Value = "Visão normal (>= 0.8)"
$ws.Cells.Item(53, 25).Value = "Não"
$ws.Cells.Item(53, 26).Value = "Um pouco"
$ws.Cells.Item(53, 27).Value = "Lenovo Vibe K5"
$ws.Cells.Item(53, 28).Value = "Android 5.1.1"
$ws.Cells.Item(53, 29).Value = "Sim"
$ws.Cells.Item(53, 30).Value = "Não preciso"
$ws.Cells.Item(53, 31).Value = "Normalmente direita"

# Row 54
$ws.Cells.Item(54, 1).Value = 42531.55948783565
$ws.Cells.Item(54, 2).Value = "Jonas Deyson Brito dos Santos"
$ws.Cells.Item(54, 3).Value = 32.0
$ws.Cells.Item(54, 4).Value = "Masculino"
$ws.Cells.Item(54, 5).Value = "Pós-graduação (Mestrado/Doutorado/Pós-doc) COMPLETO"
$ws.Cells.Item(54, 6).Value = "Computação"
$ws.Cells.Item(54, 7).Value = "Sempre direita"
$ws.Cells.Item(54, 8).Value = "Sempre direita"
$ws.Cells.Item(54, 9).Value = "Normalmente direita"
$ws.Cells.Item(54, 10).Value = "Sempre direita"
$ws.Cells.Item(54, 11).Value = "Normalmente direita"
$ws.Cells.Item(54, 12).Value = "Normalmente direita"
$ws.Cells.Item(54, 13).Value = "Sempre direita"
$ws.Cells.Item(54, 14).Value = "de 4 a 8 horas"
$ws.Cells.Item(54, 15).Value = "até 4 horas"
$ws.Cells.Item(54, 16).Value = "Nunca"
$ws.Cells.Item(54, 17).Value = "Nunca"
$ws.Cells.Item(54, 18).Value = 4.0
$ws.Cells.Item(54, 19).Value = "de 4 a 8 horas"
$ws.Cells.Item(54, 20).Value = "Nunca"
$ws.Cells.Item(54, 21).Value = "Não"
$ws.Cells.Item(54, 22).Value = 1.0
$ws.Cells.Item(54, 23).Value = "Não"
$ws.Cells.Item(54, 24).Value = "Visão normal (>= 0.8)"
$ws.Cells.Item(54, 25).Value = "Não"
$ws.Cells.Item(54, 26).Value = "Não"
$ws.Cells.Item(54, 27).Value = "Moto G 1ª geração"
$ws.Cells.Item(54, 28).Value = "Android 5.1"
$ws.Cells.Item(54, 29).Value = "Sim"
$ws.Cells.Item(54, 30).Value = "10 Jun, 13:30"
$ws.Cells.Item(54, 31).Value = "Normalmente direita"

# Row 55
$ws.Cells.Item(55, 1).Value = 42531.56238877315
$ws.Cells.Item(55, 2).Value = "Alex Reimann Cunha Lima"
$ws.Cells.Item(55, 3).Value = 33.0
$ws.Cells.Item(55, 4).Value = "Masculino"
$ws.Cells.Item(55, 5).Value = "Pós-graduação (Mestrado/Doutorado/Pós-doc) INCOMPLETO"
$ws.Cells.Item(55, 6).Value = "Direito e Ciência Computação"
$ws.Cells.Item(55, 7).Value = "Sempre direita"
$ws.Cells.Item(55, 8).Value = "Sempre direita"
$ws.Cells.Item(55, 9).Value = "Normalmente direita"
$ws.Cells.Item(55, 10).Value = "Sempre direita"
$ws.Cells.Item(55, 11).Value = "Sempre direita"
$ws.Cells.Item(55, 12).Value = "Sempre direita"
$ws.Cells.Item(55, 13).Value = "Sempre direita"
$ws.Cells.Item(55, 14).Value = "até 4 horas"
$ws.Cells.Item(55, 15).Value = "mais que 12 horas"
$ws.Cells.Item(55, 16).Value = "Nunca"
$ws.Cells.Item(55, 17).Value = "Nunca"
$ws.Cells.Item(55, 18).Value = 3.0
$ws.Cells.Item(55, 19).Value = "Nunca"
$ws.Cells.Item(55, 20).Value = "Nunca"
$ws.Cells.Item(55, 21).Value = "Wiimote"
$ws.Cells.Item(55, 22).Value = 2.0
$ws.Cells.Item(55, 23).Value = "Não"
$ws.Cells.Item(55, 24).Value = "Visão normal (>= 0.8)"
$ws.Cells.Item(55, 25).Value = "Miopia e astigmatismo"
$ws.Cells.Item(55, 26).Value = "Não"
$ws.Cells.Item(55, 27).Value = "Iphone 4, Samsung Duos e LG A275"
$ws.Cells.Item(55, 28).Value = "iOS, Android e Sistema Próprio da LG"
$ws.Cells.Item(55, 29).Value = "Sim"
$ws.Cells.Item(55, 30).Value = "OK"
$ws.Cells.Item(55, 31).Value = "Normalmente direita"

# Row 56
$ws.Cells.Item(56, 1).Value = 42531.681609050924
$ws.Cells.Item(56, 2).Value = "Mathias Fassini Mantelli"
$ws.Cells.Item(56, 3).Value = 22.0
$ws.Cells.Item(56, 4).Value = "Masculino"
$ws.Cells.Item(56, 5).Value = "Pós-graduação (Mestrado/Doutorado/Pós-doc) INCOMPLETO"
$ws.Cells.Item(56, 6).Value = "Ciência da Computação"
$ws.Cells.Item(56, 7).Value = "Sempre esquerda"
$ws.Cells.Item(56, 8).Value = "Normalmente esquerda"
$ws.Cells.Item(56, 9).Value = "Normalmente direita"
$ws.Cells.Item(56, 10).Value = "Normalmente esquerda"
$ws.Cells.Item(56, 11).Value = "Sempre esquerda"
$ws.Cells.Item(56, 12).Value = "Sempre esquerda"
$ws.Cells.Item(56, 13).Value = "Normalmente direita"
$ws.Cells.Item(56, 14).Value = "até 4 horas"
$ws.Cells.Item(56, 15).Value = "mais que 12 horas"
$ws.Cells.Item(56, 16).Value = "Nunca"
$ws.Cells.Item(56, 17).Value = "Raramente"
$ws.Cells.Item(56, 18).Value = 4.0
$ws.Cells.Item(56, 19).Value = "até 4 horas"
$ws.Cells.Item(56, 20).Value = "Nunca"
$ws.Cells.Item(56, 21).Value = "Sim, Kinect"
$ws.Cells.Item(56, 22).Value = 2.0
$ws.Cells.Item(56, 23).Value = "Não"
$ws.Cells.Item(56, 24).Value = "Visão normal (>= 0.8)"
$ws.Cells.Item(56, 25).Value = "Não"
$ws.Cells.Item(56, 26).Value = "Não"
$ws.Cells.Item(56, 27).Value = "Samsung Galaxy Note I"
$ws.Cells.Item(56, 28).Value = "Android 4.1.2"
$ws.Cells.Item(56, 29).Value = "Sim"
$ws.Cells.Item(56, 30).Value = "Ok"
$ws.Cells.Item(56, 31).Value = "Normalmente esquerda"

# Row 57
$ws.Cells.Item(57, 1).Value = 42531.6829569213
$ws.Cells.Item(57, 2).Value = "Tatiane Sequerra Stivelman"
$ws.Cells.Item(57, 3).Value = 19.0
$ws.Cells.Item(57, 4).Value = "Feminino"
$ws.Cells.Item(57, 5).Value = "Ensino superior INCOMPLETO"
$ws.Cells.Item(57, 6).Value = "Engenharia de Computação"
$ws.Cells.Item(57, 7).Value = "Sempre direita"
$ws.Cells.Item(57, 8).Value = "Sempre direita"
$ws.Cells.Item(57, 9).Value = "Sempre direita"
$ws.Cells.Item(57, 10).Value = "Sem prefêrencia"
$ws.Cells.Item(57, 11).Value = "Sempre direita"
$ws.Cells.Item(57, 12).Value = "Sempre direita"
$ws.Cells.Item(57, 13).Value = "Sempre direita"
$ws.Cells.Item(57, 14).Value = "até 4 horas"
$ws.Cells.Item(57, 15).Value = "mais que 12 horas"
$ws.Cells.Item(57, 16).Value = "de 8 a 12 horas"
$ws.Cells.Item(57, 17).Value = "Raramente"
$ws.Cells.Item(57, 18).Value = 4.0
$ws.Cells.Item(57, 19).Value = "até 4 horas"
$ws.Cells.Item(57, 20).Value = "Nunca"
$ws.Cells.Item(57, 21).Value = "Nunca"
$ws.Cells.Item(57, 22).Value = 1.0
$ws.Cells.Item(57, 23).Value = "Não"
$ws.Cells.Item(57, 24).Value = "Visão normal (>= 0.8)"
$ws.Cells.Item(57, 25).Value = "Não"
$ws.Cells.Item(57, 26).Value = "Não"
$ws.Cells.Item(57, 27).Value = "Samsung Galaxy J5"
$ws.Cells.Item(57, 28).Value = "Android 5.1.1"
$ws.Cells.Item(57, 29).Value = "Sim"
$ws.Cells.Item(57, 30).Value = "Ok"
$ws.Cells.Item(57, 31).Value = "Sem prefêrencia"

# Row 58
$ws.Cells.Item(58, 1).Value = 42531.6844531713
$ws.Cells.Item(58, 2).Value = "Fernanda Caroline Silveira Rodrigues"
$ws.Cells.Item(58, 3).Value = 27.0
$ws.Cells.Item(58, 4).Value = "Feminino"
$ws.Cells.Item(58, 5).Value = "Pós-graduação (Mestrado/Doutorado/Pós-doc) INCOMPLETO"
$ws.Cells.Item(58, 6).Value = "Ciência da Computação"
$ws.Cells.Item(58, 7).Value = "Normalmente esquerda"
$ws.Cells.Item(58, 8).Value = "Sem prefêrencia"
$ws.Cells.Item(58, 9).Value = "Normalmente esquerda"
$ws.Cells.Item(58, 10).Value = "Sempre esquerda"
$ws.Cells.Item(58, 11).Value = "Normalmente esquerda"
$ws.Cells.Item(58, 12).Value = "Sempre esquerda"
$ws.Cells.Item(58, 13).Value = "Normalmente direita"
$ws.Cells.Item(58, 14).Value = "Nunca"
$ws.Cells.Item(58, 15).Value = "de 4 a 8 horas"
$ws.Cells.Item(58, 16).Value = "Nunca"
$ws.Cells.Item(58, 17).Value = "Ocasionalmente"
$ws.Cells.Item(58, 18).Value = 2.0
$ws.Cells.Item(58, 19).Value = "Nunca"
$ws.Cells.Item(58, 20).Value = "Nunca"
$ws.Cells.Item(58, 21).Value = "Sim. Kinect, Rift."
$ws.Cells.Item(58, 22).Value = 3.0
$ws.Cells.Item(58, 23).Value = "Não."
$ws.Cells.Item(58, 24).Value = "Visão normal (>= 0.8)"
$ws.Cells.Item(58, 25).Value = "Miopia e astigmatismo."
$ws.Cells.Item(58, 26).Value = "Vertigem."
$ws.Cells.Item(58, 27).Value = "IPhone 4s"
$ws.Cells.Item(58, 28).Value = "iOs "
$ws.Cells.Item(58, 29).Value = "Sim"
$ws.Cells.Item(58, 30).Value = "Ok."
$ws.Cells.Item(58, 31).Value = "Normalmente esquerda"

# Row 59
$ws.Cells.Item(59, 1).Value = 42531.684989953705
$ws.Cells.Item(59, 2).Value = "Mariane Teixeira Giambastiani"
$ws.Cells.Item(59, 3).Value = 24.0
$ws.Cells.Item(59, 4).Value = "Feminino"
$ws.Cells.Item(59, 5).Value = "Ensino superior INCOMPLETO"
$ws.Cells.Item(59, 6).Value = "Engenharia de Computação"
$ws.Cells.Item(59, 7).Value = "Sempre direita"
$ws.Cells.Item(59, 8).Value = "Sempre direita"
$ws.Cells.Item(59, 9).Value = "Sempre direita"
$ws.Cells.Item(59, 10).Value = "Sempre direita"
$ws.Cells.Item(59, 11).Value = "Sempre direita"
$ws.Cells.Item(59, 12).Value = "Sempre direita"
$ws.Cells.Item(59, 13).Value = "Sempre direita"
$ws.Cells.Item(59, 14).Value = "Nunca"
$ws.Cells.Item(59, 15).Value = "até 4 horas"
$ws.Cells.Item(59, 16).Value = "até 4 horas"
$ws.Cells.Item(59, 17).Value = "Ocasionalmente"
$ws.Cells.Item(59, 18).Value = 2.0
$ws.Cells.Item(59, 19).Value = "Nunca"
$ws.Cells.Item(59, 20).Value = "Nunca"
$ws.Cells.Item(59, 21).Value = "wii mote"
$ws.Cells.Item(59, 22).Value = 1.0
$ws.Cells.Item(59, 23).Value = "não"
$ws.Cells.Item(59, 24).Value = "Visão normal (>= 0.8)"
$ws.Cells.Item(59, 25).Value = "Estrabismo"
$ws.Cells.Item(59, 26).Value = "não"
$ws.Cells.Item(59, 27).Value = "Sansung S3 Duo"
$ws.Cells.Item(59, 28).Value = "Android"
$ws.Cells.Item(59, 29).Value = "Sim"
$ws.Cells.Item(59, 30).Value = "ok"
$ws.Cells.Item(59, 31).Value = "Sempre direita"

# Row 60
$ws.Cells.Item(60, 1).Value = 42531.71714418981
$ws.Cells.Item(60, 2).Value = "Diego Pittol"
$ws.Cells.Item(60, 3).Value = 26.0
$ws.Cells.Item(60, 4).Value = "Masculino"
$ws.Cells.Item(60, 5).Value = "Pós-graduação (Mestrado/Doutorado/Pós-doc) INCOMPLETO"
$ws.Cells.Item(60, 6).Value = "Engenharia de Computação"
$ws.Cells.Item(60, 7).Value = "Sempre direita"
$ws.Cells.Item(60, 8).Value = "Normalmente direita"
$ws.Cells.Item(60, 9).Value = "Normalmente direita"
$ws.Cells.Item(60, 10).Value = "Normalmente direita"
$ws.Cells.Item(60, 11).Value = "Normalmente direita"
$ws.Cells.Item(60, 12).Value = "Normalmente direita"
$ws.Cells.Item(60, 13).Value = "Sempre direita"
$ws.Cells.Item(60, 14).Value = "até 4 horas"
$ws.Cells.Item(60, 15).Value = "de 4 a 8 horas"
$ws.Cells.Item(60, 16).Value = "até 4 horas"
$ws.Cells.Item(60, 17).Value = "Nunca"
$ws.Cells.Item(60, 18).Value = 3.0
$ws.Cells.Item(60, 19).Value = "Nunca"
$ws.Cells.Item(60, 20).Value = "até 4 horas"
$ws.Cells.Item(60, 21).Value = "Sim, Kinect"
$ws.Cells.Item(60, 22).Value = 1.0
$ws.Cells.Item(60, 23).Value = "Não"
$ws.Cells.Item(60, 24).Value = "Pouca perda de visão (< 0.8 e >= 0.3)"
$ws.Cells.Item(60, 25).Value = "Astigmatismo"
$ws.Cells.Item(60, 26).Value = "Não"
$ws.Cells.Item(60, 27).Value = "Moto E"
$ws.Cells.Item(60, 28).Value = "Android 5.0.0"
$ws.Cells.Item(60, 29).Value = "Sim"
$ws.Cells.Item(60, 30).Value = "Ok"
$ws.Cells.Item(60, 31).Value = "Normalmente direita"

# Row 61
$ws.Cells.Item(61, 1).Value = 42531.7479944213
$ws.Cells.Item(61, 2).Value = "Maurício Calegari Xavier"
$ws.Cells.Item(61, 3).Value = 18.0
$ws.Cells.Item(61, 4).Value = "Masculino"
$ws.Cells.Item(61, 5).Value = "Ensino superior INCOMPLETO"
$ws.Cells.Item(61, 6).Value = "Engenharia de Computação"
$ws.Cells.Item(61, 7).Value = "Sempre direita"
$ws.Cells.Item(61, 8).Value = "Sempre direita"
$ws.Cells.Item(61, 9).Value = "Sempre direita"
$ws.Cells.Item(61, 10).Value = "Sempre direita"
$ws.Cells.Item(61, 11).Value = "Normalmente direita"
$ws.Cells.Item(61, 12).Value = "Normalmente direita"
$ws.Cells.Item(61, 13).Value = "Sempre direita"
$ws.Cells.Item(61, 14).Value = "Nunca"
$ws.Cells.Item(61, 15).Value = "de 4 a 8 horas"
$ws.Cells.Item(61, 16).Value = "até 4 horas"
$ws.Cells.Item(61, 17).Value = "Raramente"
$ws.Cells.Item(61, 18).Value = 4.0
$ws.Cells.Item(61, 19).Value = "Nunca"
$ws.Cells.Item(61, 20).Value = "Nunca"
$ws.Cells.Item(61, 21).Value = "Sim. Kinect."
$ws.Cells.Item(61, 22).Value = 3.0
$ws.Cells.Item(61, 23).Value = "Não"
$ws.Cells.Item(61, 24).Value = "Pouca perda de visão (< 0.8 e >= 0.3)"
$ws.Cells.Item(61, 25).Value = "Astigmatismo."
$ws.Cells.Item(61, 26).Value = "Não."
$ws.Cells.Item(61, 27).Value = "Samsung Galaxy Win"
$ws.Cells.Item(61, 28).Value = "Android 4.2"
$ws.Cells.Item(61, 29).Value = "Sim"
$ws.Cells.Item(61, 30).Value = "Ok."
$ws.Cells.Item(61, 31).Value = "Sempre direita"

# Row 62
$ws.Cells.Item(62, 1).Value = 42531.74828921296
$ws.Cells.Item(62, 2).Value = "Emanuel Teribele Novakoski"
$ws.Cells.Item(62, 3).Value = 18.0
$ws.Cells.Item(62, 4).Value = "Masculino"
$ws.Cells.Item(62, 5).Value = "Ensino superior INCOMPLETO"
$ws.Cells.Item(62, 6).Value = "Engenharia de Computacao"
$ws.Cells.Item(62, 7).Value = "Sempre direita"
$ws.Cells.Item(62, 8).Value = "Sempre direita"
$ws.Cells.Item(62, 9).Value = "Sempre direita"
$ws.Cells.Item(62, 10).Value = "Sempre direita"
$ws.Cells.Item(62, 11).Value = "Sempre direita"
$ws.Cells.Item(62, 12).Value = "Sempre direita"
$ws.Cells.Item(62, 13).Value = "Sempre direita"
$ws.Cells.Item(62, 14).Value = "mais que 12 horas"
$ws.Cells.Item(62, 15).Value = "de 4 a 8 horas"
$ws.Cells.Item(62, 16).Value = "até 4 horas"
$ws.Cells.Item(62, 17).Value = "Raramente"
$ws.Cells.Item(62, 18).Value = 4.0
$ws.Cells.Item(62, 19).Value = "de 4 a 8 horas"
$ws.Cells.Item(62, 20).Value = "Nunca"
$ws.Cells.Item(62, 21).Value = "Oculus Rift"
$ws.Cells.Item(62, 22).Value = 1.0
$ws.Cells.Item(62, 23).Value = "Nao"
$ws.Cells.Item(62, 24).Value = "Perda de visão moderada (< 0.3 e >= 0.125)"
$ws.Cells.Item(62, 25).Value = "Astigmatismo"
$ws.Cells.Item(62, 26).Value = "Nao"
$ws.Cells.Item(62, 27).Value = "LG G3"
$ws.Cells.Item(62, 28).Value = "Android 4.4"
$ws.Cells.Item(62, 29).Value = "Sim"
$ws.Cells.Item(62, 30).Value = "ok"
$ws.Cells.Item(62, 31).Value = "Normalmente direita"

# Row 63
$ws.Cells.Item(63, 1).Value = 42531.748840775464
$ws.Cells.Item(63, 2).Value = "Mauricio Barbosa da Rocha"
$ws.Cells.Item(63, 3).Value = 24.0
$ws.Cells.Item(63, 4).Value = "Masculino"
$ws.Cells.Item(63, 5).Value = "Ensino superior INCOMPLETO"
$ws.Cells.Item(63, 6).Value = "Engenharia de computacao"
$ws.Cells.Item(63, 7).Value = "Sempre direita"
$ws.Cells.Item(63, 8).Value = "Sempre direita"
$ws.Cells.Item(63, 9).Value = "Sempre direita"
$ws.Cells.Item(63, 10).Value = "Sempre direita"
$ws.Cells.Item(63, 11).Value = "Sempre direita"
$ws.Cells.Item(63, 12).Value = "Sempre direita"
$ws.Cells.Item(63, 13).Value = "Normalmente esquerda"
$ws.Cells.Item(63, 14).Value = "até 4 horas"
$ws.Cells.Item(63, 15).Value = "de 4 a 8 horas"
$ws.Cells.Item(63, 16).Value = "Nunca"
$ws.Cells.Item(63, 17).Value = "Nunca"
$ws.Cells.Item(63, 18).Value = 2.0
$ws.Cells.Item(63, 19).Value = "Nunca"
$ws.Cells.Item(63, 20).Value = "Nunca"
$ws.Cells.Item(63, 21).Value = "Sim"
$ws.Cells.Item(63, 22).Value = 1.0
$ws.Cells.Item(63, 23).Value = "Nao"
$ws.Cells.Item(63, 24).Value = "Visão normal (>= 0.8)"
$ws.Cells.Item(63, 25).Value = "Nao"
$ws.Cells.Item(63, 26).Value = "Nao"
$ws.Cells.Item(63, 27).Value = "samsung"
$ws.Cells.Item(63, 28).Value = "android"
$ws.Cells.Item(63, 29).Value = "Sim"
$ws.Cells.Item(63, 30).Value = "ok"
$ws.Cells.Item(63, 31).Value = "Sempre direita"

# Row 64
$ws.Cells.Item(64, 1).Value = 42533.631650613424
$ws.Cells.Item(64, 2).Value = "Guilherme Fonseca Ribeiro"
$ws.Cells.Item(64, 3).Value = 24.0
$ws.Cells.Item(64, 4).Value = "Masculino"
$ws.Cells.Item(64, 5).Value = "Ensino superior INCOMPLETO"
$ws.Cells.Item(64, 6).Value = "Ciência da Computação"
$ws.Cells.Item(64, 7).Value = "Sempre direita"
$ws.Cells.Item(64, 8).Value = "Sempre direita"
$ws.Cells.Item(64, 9).Value = "Sempre direita"
$ws.Cells.Item(64, 10).Value = "Sempre direita"
$ws.Cells.Item(64, 11).Value = "Sempre direita"
$ws.Cells.Item(64, 12).Value = "Sempre direita"
$ws.Cells.Item(64, 13).Value = "Sempre direita"
$ws.Cells.Item(64, 14).Value = "Nunca"
$ws.Cells.Item(64, 15).Value = "mais que 12 horas"
$ws.Cells.Item(64, 16).Value = "Nunca"
$ws.Cells.Item(64, 17).Value = "Nunca"
$ws.Cells.Item(64, 18).Value = 4.0
$ws.Cells.Item(64, 19).Value = "Nunca"
$ws.Cells.Item(64, 20).Value = "Nunca"
$ws.Cells.Item(64, 21).Value = "Nunca"
$ws.Cells.Item(64, 22).Value = 1.0
$ws.Cells.Item(64, 23).Value = "Não"
$ws.Cells.Item(64, 24).Value = "Pouca perda de visão (< 0.8 e >= 0.3)"
$ws.Cells.Item(64, 25).Value = "Astigmatismo e miopia"
$ws.Cells.Item(64, 26).Value = "Não"
$ws.Cells.Item(64, 27).Value = "iPhone 4s"
$ws.Cells.Item(64, 28).Value = "Alguma versão do iOs"
$ws.Cells.Item(64, 29).Value = "Sim"
$ws.Cells.Item(64, 30).Value = "Ok"
$ws.Cells.Item(64, 31).Value = "Normalmente direita"
